# ---------------------------------------------------------------------------
# [Item list] Create header added. Not tested
#
# 1. Rename "Sheet1" -> "Funkcje"
# 2. Add new sheet "Arch" (sheetId 3) right after "Funkcje"
# 3. Populate "Arch" with the Item-list / architecture table (header row + data)
# 4. Center the header row of "Arch"
# 5. Give row 14 on "Funkcje" extra height (wrapped description got taller)
# 6. Update stored selections / active sheet / zoom to match the new layout
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$tabelle1 = $wb.Worksheets.Item(1)
$funkcje  = $wb.Worksheets.Item(2)

# --- rename Sheet1 -> Funkcje ------------------------------------------------
$funkcje.Name = "Funkcje"

# --- add the new "Arch" sheet right after "Funkcje" -------------------------
$arch = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $funkcje)
$arch.Name = "Arch"

# --- "Funkcje": row 14 grew taller (wrapped text needs two lines now) -------
$funkcje.Rows.Item(14).RowHeight = 30

# --- "Arch": header row (row 3) ---------------------------------------------
$arch.Range("B3").Value = "App code"
$arch.Range("C3").Value = "App prog interface"
$arch.Range("D3").Value = "Middleware"
$arch.Range("F3").Value = "HAL"
$arch.Range("B3:D3").HorizontalAlignment = -4108
$arch.Range("F3").HorizontalAlignment = -4108

# --- "Arch": data rows 4-17 --------------------------------------------------
$arch.Range("D4").Value = "void LCD_Init_HW(...)"
$arch.Range("E4").Value = "void MX_SPI2_Init(...)"
$arch.Range("F4").Value = "void HAL_StatusTypeDef HAL_SPI_Transmit(...)"

$arch.Range("D5").Value = "void LCD_Configure(...)"
$arch.Range("E5").Value = "void SPI_Send_Command(...)"

$arch.Range("D6").Value = "void Set_Address (...)"
$arch.Range("E6").Value = "void SPI_Send_Data_8bit(...)"

$arch.Range("D7").Value = "void LCD_Data_Preparation(...)"
$arch.Range("E7").Value = "void SPI_Send_Data_16bit(...)"

$arch.Range("D8").Value  = "void Fill_display(...)"
$arch.Range("D9").Value  = "void LCD_Init(...)"
$arch.Range("D10").Value = "void Draw_Point(...)"
$arch.Range("D11").Value = "void LCD_DrawLine(...)"
$arch.Range("D12").Value = "void LCD_DrawRectangle(...)"
$arch.Range("D13").Value = "void LCD_DrawCircle(...)"
$arch.Range("D14").Value = "void LCD_DisplayChar(…)"
$arch.Range("D15").Value = "void LCD_DisplayString (...)"
$arch.Range("D16").Value = "void LCD_DisplayNum(...)"
$arch.Range("D17").Value = "void LCD_Features_Selftest(...)"

# --- "Arch": column widths (characters) -------------------------------------
$arch.Columns.Item(2).ColumnWidth = 16.666666666666668
$arch.Columns.Item(3).ColumnWidth = 28.666666666666668
$arch.Columns.Item(4).ColumnWidth = 27.333333333333332
$arch.Columns.Item(5).ColumnWidth = 43

# --- selections / active cells ----------------------------------------------
$tabelle1.Range("E29").Select() | Out-Null
$funkcje.Range("D19").Select() | Out-Null
$arch.Range("D19").Select() | Out-Null

# --- zoom + activate the new sheet ------------------------------------------
$arch.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 85

Write-Output "done"
